$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '68.873.00'
$ws.Range("E2").Value = '  +1.44%  '

$ws.Range("D3").Value = '3.282.89'
$ws.Range("E3").Value = '  +0.50%  '

$ws.Range("E4").Value = '  +0.01%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '583.29'
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = '  +0.47%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '186.11'
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = '  +2.38%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.599'
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = '  -0.34%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.133'
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = '  -0.46%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '6.66'
$ws.Range("D10").ClearFormats()

$ws.Range("E11").Value = '  +1.44%  '

$ws.Range("D12").Value = '3.850.13'
$ws.Range("E12").Value = '  +0.48%  '

$ws.Range("E13").Value = '  -0.13%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '28.50'
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = '  +0.20%  '

$ws.Range("D15").Value = '68.772.31'
$ws.Range("E15").Value = '  +1.42%  '

$ws.Range("E16").Value = '  +1.69%  '

$ws.Range("D17").Value = '3.255.90'
$ws.Range("E17").Value = '  -0.70%  '

$ws.Range("E18").Value = '  +0.27%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '13.63'
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = '  +0.89%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '393.34'
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = '  +4.57%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '7.74'
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = '  +1.20%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '71.87'
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = '  +0.80%  '

$ws.Range("E23").Value = '  +0.04%  '

$ws.Range("E24").Value = '  +1.47%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.0000121'
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = '  +0.33%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.189'
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = '  +4.77%  '

$ws.Range("E27").Value = '  +0.72%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '0.998'
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = '  +0.04%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '1.99'
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = '  +0.36%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '5.75'
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = '  +1.35%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '23.15'
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = '  +1.69%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '7.19'
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = '  +3.93%  '

$ws.Range("E33").Value = '  +2.91%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.998'
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = '  +0.04%  '

$ws.Range("B35").Value = 'ImmutableX'
$ws.Range("C35").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.52'
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = '  +0.16%  '

$ws.Range("B36").Value = 'Monero'
$ws.Range("C36").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '163.93'
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = '  +1.20%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '1.96'
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = '  +6.29%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.827'
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = '  -2.83%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '26.82'
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = '  +0.23%  '

$ws.Range("E40").Value = '  -0.77%  '

$ws.Range("E41").Value = '  -2.63%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '2.57'
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = '  -1.85%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '41.46'
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = '  +1.59%  '

$ws.Range("B44").Value = 'InjectiveProtocol'
$ws.Range("C44").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '25.52'
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = '  -0.76%  '

$ws.Range("B45").Value = 'Hedera'
$ws.Range("C45").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.0691'
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = '  +1.87%  '

$ws.Range("B46").Value = 'Maker'
$ws.Range("C46").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D46").Value = '2.655.18'
$ws.Range("E46").Value = '  -1.44%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '342.46'
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = '  -2.43%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.0283'
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = '  +0.78%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '32.10'
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = '  +2.67%  '

$ws.Range("E50").Value = '  +3.28%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.998'
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = '  -0.23%  '
